$d = $word.ActiveDocument
# Step 1: delete paragraph 11 ("Sin datos") entirely
$d.Paragraphs(11).Range.Delete()

# Step 2: paragraph 10 ("Datos:") -> Precondicion paragraph + new empty paragraph
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Sinespaciado"/><w:tabs><w:tab w:val="left" w:pos="4345"/></w:tabs><w:rPr><w:b/><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="es-AR"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Precondición: </w:t></w:r><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Se debe haber ejecutado con éxito el CP0001, CP0004, CP0006 y CP0013.</w:t></w:r><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:tab/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Sinespaciado"/><w:tabs><w:tab w:val="left" w:pos="4345"/></w:tabs><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(10).Range.InsertXML($xml2)

# Step 3: paragraph 9 ("Precondicion...") -> Objetivo paragraph + new empty paragraph
$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Sinespaciado"/><w:tabs><w:tab w:val="left" w:pos="7025"/></w:tabs><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Objetivo: </w:t></w:r><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Esperar que el sistema cree una solicitud de amistad para el usuario que se agrega.</w:t></w:r><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:tab/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Sinespaciado"/><w:tabs><w:tab w:val="left" w:pos="7025"/></w:tabs><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(9).Range.InsertXML($xml3)

# Step 4: paragraph 8 ("Datos de prueba: ") -> CP0014 (Ttulo4) + new empty paragraph
$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Ttulo4"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>CP0014</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(8).Range.InsertXML($xml4)

# Step 5: paragraph 6 ("Datos de prueba") -> split into "Cas"+"os de prueba" runs
$xml5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Ttulo3"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Cas</w:t></w:r><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>os de prueba</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(6).Range.InsertXML($xml5)

# Step 6: paragraph 3 (empty Textbody) -> add lang rPr to pPr
$xml6 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Textbody"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(3).Range.InsertXML($xml6)

# Step 7: paragraph 2 ("Procedimiento de prueba") -> add lang rPr to pPr
$xml7 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Ttulo3"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Procedimiento de prueba</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(2).Range.InsertXML($xml7)

Write-Host "Done. Paragraph count: " $d.Paragraphs.Count
